# Updated utility and disutility data
# This script reproduces the data edits made to the "ae_disutility" sheet:
#  - C column: mean disutility values become explicit formulas (difference of two utilities)
#  - D column: se values become formulas (ROUND(0.33*mean,3)) i.e. se assumed 1/3rd of mean
#  - E column: reference updated from nafees2008health to nafees2017health (except row 8,
#    which now cites doyle2008health)
#  - F column: notes text updated, and the note in row 8 is highlighted in red font
#  - Selection/active sheet changed back to state_utility (sheet1, cell B3)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # state_utility
$ws2 = $wb.Worksheets.Item(2)   # ae_disutility

# ---------------------------------------------------------------------------
# ae_disutility (sheet2) - column C: mean = formula (utility at baseline - utility w/ AE)
# ---------------------------------------------------------------------------
$ws2.Range("C2").Formula  = "=0.754-0.46"
$ws2.Range("C3").Formula  = "=0.754-0.46"
$ws2.Range("C4").Formula  = "=0.754-0.532"
$ws2.Range("C5").Formula  = "=0.754-0.603"
$ws2.Range("C6").Formula  = "=0.754-0.616"
$ws2.Range("C7").Formula  = "=0.754-0.603"
$ws2.Range("C8").Formula  = "=0.626-0.461"
$ws2.Range("C9").Formula  = "=0.754-0.603"
$ws2.Range("C10").Formula = "=0.754-0.603"
$ws2.Range("C11").Formula = "=0.754-0.603"

# ---------------------------------------------------------------------------
# ae_disutility (sheet2) - column D: se = ROUND(0.33*mean,3)  ("se assumed 1/3rd of mean")
# D2 keeps its own formula; D3:D11 are entered together so they share one formula group,
# exactly like the original authoring in Excel.
# ---------------------------------------------------------------------------
$ws2.Range("D2").Formula   = "=ROUND(0.33*C2,3)"
$ws2.Range("D3:D11").Formula = "=ROUND(0.33*C3,3)"

# ---------------------------------------------------------------------------
# ae_disutility (sheet2) - column E: reference updated
# ---------------------------------------------------------------------------
$ws2.Range("E2").Value  = "nafees2017health"
$ws2.Range("E3").Value  = "nafees2017health"
$ws2.Range("E4").Value  = "nafees2017health"
$ws2.Range("E5").Value  = "nafees2017health"
$ws2.Range("E6").Value  = "nafees2017health"
$ws2.Range("E7").Value  = "nafees2017health"
$ws2.Range("E8").Value  = "doyle2008health"
$ws2.Range("E9").Value  = "nafees2017health"
$ws2.Range("E10").Value = "nafees2017health"
$ws2.Range("E11").Value = "nafees2017health"

# ---------------------------------------------------------------------------
# ae_disutility (sheet2) - column F: notes updated
# ---------------------------------------------------------------------------
$ws2.Range("F2").Value  = "Assumed equal to fatigue disutility; se assumed 1/3rd of mean"
$ws2.Range("F3").Value  = "Assumed equal to fatigue disutility; se assumed 1/3rd of mean"
$ws2.Range("F5").Value  = "Assumed equal to rash disutility; se assumed 1/3rd of mean"
$ws2.Range("F6").Value  = "Assumed equal to hair loss; se assumed 1/3rd of mean"
$ws2.Range("F7").Value  = "Assumed equal to rash disutility; se assumed 1/3rd of mean"
$ws2.Range("F8").Value  = "Assumed equal to cough, dyspnea, and pain; se assumed 1/3rd of mean"
$ws2.Range("F9").Value  = "Assumed equal to rash disutility; se assumed 1/3rd of mean"
$ws2.Range("F11").Value = "Assumed equal to rash disutility; se assumed 1/3rd of mean"

# The note in F8 is now flagged/highlighted with a red font color.
$ws2.Range("F8").Font.Color = 255

# ---------------------------------------------------------------------------
# Restore selection: state_utility becomes the active sheet again, with B3 selected;
# ae_disutility keeps D2 selected (and is no longer the active/visible tab).
# ---------------------------------------------------------------------------
$ws2.Range("D2").Select()
$ws1.Activate()
$ws1.Range("B3").Select()
